$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns: "location" -> "name", "network" -> "netid"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "netid"

# Update the active selection from E2 to C1
$ws.Range("C1").Select()
